$wb = $excel.ActiveWorkbook

$wsGlobals      = $wb.Worksheets.Item("Globals")
$wsPowerPlants  = $wb.Worksheets.Item("PowerPlants")
$wsFuels        = $wb.Worksheets.Item("Fuels")
$wsConnections  = $wb.Worksheets.Item("Connections")

# --- Globals sheet: DiscountRate is no longer included ---
$wsGlobals.Range("B2").Value = "N"

# --- PowerPlants sheet: several variables no longer included ---
# (ExpectedLifetime, row 3, stays "Y" / always included since it will always be integer)
$wsPowerPlants.Range("B2").Value = "N"
$wsPowerPlants.Range("B4").Value = "N"
$wsPowerPlants.Range("B5").Value = "N"
$wsPowerPlants.Range("B6").Value = "N"
$wsPowerPlants.Range("B7").Value = "N"
$wsPowerPlants.Range("B8").Value = "N"

# Highlight MaxCapacity and MaxActivity rows in yellow
$wsPowerPlants.Range("B9:B10").Interior.Color = 65535

# --- Fuels sheet: variables no longer included ---
$wsFuels.Range("B2").Value = "N"
$wsFuels.Range("B3").Value = "N"
$wsFuels.Range("B4").Value = "N"
$wsFuels.Range("B5").Value = "N"
$wsFuels.Range("B6").Value = "N"

# --- Connections sheet: variables no longer included ---
$wsConnections.Range("B2").Value = "N"
$wsConnections.Range("B3").Value = "N"
$wsConnections.Range("B4").Value = "N"
$wsConnections.Range("B5").Value = "N"

# --- Update each sheet's selection/active cell ---
$wsPowerPlants.Activate()
$wsPowerPlants.Range("B9").Select()

$wsFuels.Activate()
$wsFuels.Range("B6").Select()

$wsConnections.Activate()
$wsConnections.Range("B6").Select()

# Globals becomes the active sheet, with B3 selected
$wsGlobals.Activate()
$wsGlobals.Range("B3").Select()
